$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"
$linkA = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/09fb7839915c6404d4984363e6cf3b07567038a1/e2e/a.md"
$linkB = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/09fb7839915c6404d4984363e6cf3b07567038a1/e2e/b.md"

# ---------------------------------------------------------------------------
# Overview sheet: the per-locale "Status" columns (E = zh-cn, F = de-de)
# now report the handback is in sync with en-US.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------------
# zh-cn sheet: target file + handback file/datetime now populated for the
# handed-back "a.md" file.
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

$wsZh.Range("I2").Value = "a.md"
$wsZh.Range("I3").Value = "a.md"
$wsZh.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZh.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-30 00:38:52"
$wsZh.Range("K3").Value = "2016-08-30 00:38:52"

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $linkA, [Type]::Missing, [Type]::Missing, "a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $linkA, [Type]::Missing, [Type]::Missing, "a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $linkB, [Type]::Missing, [Type]::Missing, "b.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $linkA, [Type]::Missing, [Type]::Missing, "a.md")

$wsZh.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZh.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------------
# de-de sheet: same shape of change as zh-cn, but with its own handback
# file/datetime values.
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

$wsDe.Range("I2").Value = "a.md"
$wsDe.Range("I3").Value = "a.md"
$wsDe.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDe.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-30 00:38:59"
$wsDe.Range("K3").Value = "2016-08-30 00:38:59"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $linkA, [Type]::Missing, [Type]::Missing, "a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $linkA, [Type]::Missing, [Type]::Missing, "a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $linkB, [Type]::Missing, [Type]::Missing, "b.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $linkA, [Type]::Missing, [Type]::Missing, "a.md")

$wsDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDe.Columns.Item(10).ColumnWidth = 39.166666666666664

Write-Host "Handback report generated"
